$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove a few leftover lambda parameter values that no longer apply ---
$ws.Range("W3").ClearContents()
$ws.Range("T5").ClearContents()
$ws.Range("T6").ClearContents()

# --- Add new lambda parameter values for rows 26, 27 and 32 ---
# Row 26: Division of Water Quality (SWRCB)
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 2
$ws.Range("F26").Value = 1

# Row 27: Groundwater Management (SWRCB)
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 1

# Row 32: central valley water board
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 0.5
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = 0.5

# --- Update the selected/active cell on the sheet (view state) ---
$ws.Range("T6").Select()
